$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 382.44446
$ws.Range("I6").Value = 270.2857
$ws.Range("K6").Value = 810.8571000000001
$ws.Range("M6").Value = -698.8571000000001
$ws.Range("H8").Value = 50.666668
$ws.Range("I8").Value = 50.666668
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 152.000004
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -13.00000399999999
$ws.Range("N8").ClearContents()
$ws.Range("H20").Value = 1471
$ws.Range("I20").Value = 1471
$ws.Range("K20").Value = 1471
$ws.Range("M20").Value = -1241
$ws.Range("H34").Value = 6097.7144
$ws.Range("I34").Value = 1280.6666
$ws.Range("J34").Value = 35000
$ws.Range("K34").Value = 1280.6666
$ws.Range("L34").Value = 35000
$ws.Range("M34").Value = -1077.6666
$ws.Range("N34").Value = -35406
$ws.Range("H35").Value = 1471
$ws.Range("I35").Value = 1471
$ws.Range("K35").Value = 1471
$ws.Range("M35").Value = -1092
$ws.Range("H36").Value = 6097.7144
$ws.Range("I36").Value = 1280.6666
$ws.Range("J36").Value = 35000
$ws.Range("K36").Value = 1280.6666
$ws.Range("L36").Value = 35000
$ws.Range("M36").Value = -565.6666
$ws.Range("N36").Value = -36430
$ws.Range("H39").Value = 275.92593
$ws.Range("I39").Value = 122.14286
$ws.Range("J39").Value = 441.53845
$ws.Range("K39").Value = 366.42858
$ws.Range("L39").Value = 1324.61535
$ws.Range("M39").Value = -70.42858000000001
$ws.Range("N39").Value = -1916.61535
$ws.Range("H47").Value = 27740
$ws.Range("I47").Value = 19500
$ws.Range("J47").Value = 30486.666
$ws.Range("K47").Value = 19500
$ws.Range("L47").Value = 30486.666
$ws.Range("M47").Value = -18528
$ws.Range("N47").Value = -32430.666
$ws.Range("H54").Value = 15038
$ws.Range("I54").Value = 15038
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 15038
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -14552
$ws.Range("N54").ClearContents()
$ws.Range("H63").Value = 31208.875
$ws.Range("J63").Value = 31208.875
$ws.Range("L63").Value = 31208.875
$ws.Range("N63").Value = -32456.875
$ws.Range("H66").Value = 31208.875
$ws.Range("J66").Value = 31208.875
$ws.Range("L66").Value = 93626.625
$ws.Range("N66").Value = -99866.625
$ws.Range("H94").Value = 7392.3076
$ws.Range("I94").Value = 4075
$ws.Range("K94").Value = 4075
$ws.Range("M94").Value = -3624
$ws.Range("H100").Value = 31251044
$ws.Range("I100").Value = 1171
$ws.Range("J100").Value = 83334160
$ws.Range("K100").Value = 1171
$ws.Range("L100").Value = 83334160
$ws.Range("M100").Value = -630
$ws.Range("N100").Value = -83335242
$ws.Range("H138").Value = 2005.1719
$ws.Range("I138").Value = 1599.9615
$ws.Range("J138").Value = 2282.4211
$ws.Range("K138").Value = 4799.8845
$ws.Range("L138").Value = 6847.263300000001
$ws.Range("M138").Value = 340.1154999999999
$ws.Range("N138").Value = -17127.2633

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1745.3334
$ws.Range("I61").Value = 2098.8
$ws.Range("J61").Value = 1492.8572
$ws.Range("K61").Value = 2098.8
$ws.Range("L61").Value = 1492.8572
$ws.Range("M61").Value = -1886.8
$ws.Range("N61").Value = -1916.8572
$ws.Range("H136").Value = 1745.3334
$ws.Range("I136").Value = 2098.8
$ws.Range("J136").Value = 1492.8572
$ws.Range("K136").Value = 6296.400000000001
$ws.Range("L136").Value = 4478.571599999999
$ws.Range("M136").Value = -3746.400000000001
$ws.Range("N136").Value = -9578.571599999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5625
$ws.Range("I5").Value = 900
$ws.Range("K5").Value = 900
$ws.Range("M5").Value = -787
$ws.Range("H105").Value = 2566
$ws.Range("I105").Value = 1757.4584
$ws.Range("J105").Value = 4330.091
$ws.Range("K105").Value = 1757.4584
$ws.Range("L105").Value = 4330.091
$ws.Range("M105").Value = -10.45839999999998
$ws.Range("N105").Value = -7824.091

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H31").Value = 2169.4307
$ws.Range("I31").Value = 1635.5938
$ws.Range("J31").Value = 2687.0908
$ws.Range("K31").Value = 1635.5938
$ws.Range("L31").Value = 2687.0908
$ws.Range("M31").Value = -1340.5938
$ws.Range("N31").Value = -3277.0908
$ws.Range("H34").Value = 2169.4307
$ws.Range("I34").Value = 1635.5938
$ws.Range("J34").Value = 2687.0908
$ws.Range("K34").Value = 1635.5938
$ws.Range("L34").Value = 2687.0908
$ws.Range("M34").Value = -1433.5938
$ws.Range("N34").Value = -3091.0908

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 251.95238
$ws.Range("I26").Value = 227.28572
$ws.Range("J26").Value = 264.2857
$ws.Range("K26").Value = 681.85716
$ws.Range("L26").Value = 792.8571000000001
$ws.Range("M26").Value = -393.85716
$ws.Range("N26").Value = -1368.8571
$ws.Range("H139").Value = 1900.909
$ws.Range("I139").Value = 1322.3684
$ws.Range("J139").Value = 3194.1177
$ws.Range("K139").Value = 3967.1052
$ws.Range("L139").Value = 9582.3531
$ws.Range("M139").Value = 1172.8948
$ws.Range("N139").Value = -19862.3531

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 8000
$ws.Range("J26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("N26").Value = -8560
$ws.Range("H50").Value = 8000
$ws.Range("J50").Value = 8000
$ws.Range("L50").Value = 8000
$ws.Range("N50").Value = -8996

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 12700
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 25000
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = -287
$ws.Range("N4").Value = -25226
$ws.Range("H28").Value = 12700
$ws.Range("I28").Value = 400
$ws.Range("J28").Value = 25000
$ws.Range("K28").Value = 400
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = -168
$ws.Range("N28").Value = -25464
$ws.Range("H30").Value = 23154
$ws.Range("I30").Value = 23154
$ws.Range("K30").Value = 23154
$ws.Range("M30").Value = -23046
$ws.Range("H33").Value = 3740
$ws.Range("I33").Value = 2500
$ws.Range("J33").Value = 4566.6665
$ws.Range("K33").Value = 2500
$ws.Range("L33").Value = 4566.6665
$ws.Range("M33").Value = -2210
$ws.Range("N33").Value = -5146.6665
$ws.Range("H35").Value = 1866.6666
$ws.Range("I35").Value = 1866.6666
$ws.Range("K35").Value = 1866.6666
$ws.Range("M35").Value = -1530.6666
$ws.Range("H37").Value = 12700
$ws.Range("I37").Value = 400
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 400
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = -293
$ws.Range("N37").Value = -25214
$ws.Range("H56").Value = 20051
$ws.Range("I56").Value = 20051
$ws.Range("K56").Value = 20051
$ws.Range("M56").Value = -19360

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 7942.5
$ws.Range("I58").Value = 7885
$ws.Range("K58").Value = 7885
$ws.Range("M58").Value = -7577
$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 8000
$ws.Range("K61").Value = 8000
$ws.Range("M61").Value = -7708
$ws.Range("H100").Value = 501041.06
$ws.Range("I100").Value = 1121.6471
$ws.Range("K100").Value = 2243.2942
$ws.Range("M100").Value = -1702.2942

Write-Host "All edits applied"